# Apply the "edited ref tables for yll" change:
#  - Populate column E (LE2016 reference values) for rows 2-41 on the
#    yll_le_5yr sheet.
#  - Make yll_le_5yr the active/selected sheet (tabSelected) with the
#    new E2:E41 selection, and make yll_le_1yr no longer the selected tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("yll_le_5yr")
$ws2 = $wb.Worksheets.Item("yll_le_1yr")

# New LE2016 reference values for E2:E41 (column E), rows 2-41.
$values = @(
    83.08,
    80.878827999999999,
    76.459946000000002,
    71.505289000000005,
    66.434832,
    61.488934999999998,
    56.616292000000001,
    51.724936999999997,
    46.916581999999998,
    41.956308,
    37.257553000000001,
    32.619971999999997,
    28.101426,
    23.678177000000002,
    19.365911000000001,
    15.514853,
    11.792838,
    8.6029809999999998,
    6.0185110000000002,
    3.7965770000000001,
    79.459999999999994,
    77.300284000000005,
    72.877994999999999,
    67.925471999999999,
    62.865049999999997,
    58.021082999999997,
    53.213501999999998,
    48.387422999999998,
    43.655194000000002,
    38.802669999999999,
    34.210920000000002,
    29.710315000000001,
    25.336711000000001,
    21.12154,
    17.085767000000001,
    13.550349000000001,
    10.199809999999999,
    7.4235429999999996,
    5.2128209999999999,
    3.4080970000000002
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 5).Value = $values[$i]
}

# yll_le_5yr becomes the active sheet/tab, with E2:E41 selected.
# (yll_le_1yr's own selection stays on K13 - its existing state - it just
# stops being the tabbed/active sheet once yll_le_5yr is activated below.)
$ws1.Activate() | Out-Null
$ws1.Range("E2:E41").Select() | Out-Null
